# Update project, add new function
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Cylinder")

# Insert a new column F for "verificationDate" (shifts old brand column to G)
$ws.Range("F1").EntireColumn.Insert()

# Match the new column's width to column E so they form one contiguous
# width band (E:F), mirroring the rest of the table's column widths.
$ws.Columns("F").ColumnWidth = $ws.Columns("E").ColumnWidth

# Header for the new column
$ws.Range("F1").Value = "verificationDate"

# Update the test codes in column A (rows 2-4)
$ws.Range("A2").Value = "TEST1608004"
$ws.Range("A3").Value = "TEST1608005"
$ws.Range("A4").Value = "TEST1608006"

# Fix capitalization of "1 van" -> "1 Van" (now in column D, row 4)
$ws.Range("D4").Value = "1 Van"

# Apply the date number format to F2 first, then propagate it to F3:F4 via
# copy/paste-special so all three cells share one cell style record.
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()
$ws.Range("F3:F4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the verification dates for the new column
$ws.Range("F2").Value = Get-Date -Year 2020 -Month 8 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("F3").Value = Get-Date -Year 2020 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("F4").Value = Get-Date -Year 2020 -Month 9 -Day 4 -Hour 0 -Minute 0 -Second 0

# Update the selection on the sheet
$ws.Range("A2:A4").Select()

Write-Output "done"
